$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 30-47: fill in Description (C) / Originates from (D) columns with the
# new custom-error-code rows (Kalman filter / pose estimator related errors).
# Cells are written row-by-row, column C then D, so that new shared-string
# entries get appended in the same order they first appear.
# ---------------------------------------------------------------------------

# Row 30
$ws.Range("C30").Value = "Trouble creating new Kalman FIlter"
$ws.Range("D30").Value = "KalmanFilter_New.vi"

# Row 31
$ws.Range("C31").Value = "Trouble creating new Unscented Kalman Filter"
$ws.Range("D31").Value = "UnscentedKalmanFilter_New.vi"

# Row 32
$ws.Range("C32").Value = "Trouble creating new Linear Quadratic Regulator."
$ws.Range("D32").Value = "LinearQuadraticRegulator_New.vi"

# Row 33 (style changes from wrap/general to no-wrap/general -> copy format from D7 which already has that style)
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C33").PasteSpecial(-4122) | Out-Null
$ws.Range("C33").Value = "MeanX function reference is invalid."
$ws.Range("D33").Value = "UnscentedKalmanFilter_New_FuncGroup.vi"

# Row 34
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C34").PasteSpecial(-4122) | Out-Null
$ws.Range("C34").Value = "MeanY function reference is invalid."
$ws.Range("D34").Value = "UnscentedKalmanFilter_New_FuncGroup.vi"

# Row 35
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C35").PasteSpecial(-4122) | Out-Null
$ws.Range("C35").Value = "ResidX function reference is invalid."
$ws.Range("D35").Value = "UnscentedKalmanFilter_New_FuncGroup.vi"

# Row 36
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C36").PasteSpecial(-4122) | Out-Null
$ws.Range("C36").Value = "ResidY function reference is invalid."
$ws.Range("D36").Value = "UnscentedKalmanFilter_New_FuncGroup.vi"

# Row 37
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C37").PasteSpecial(-4122) | Out-Null
$ws.Range("C37").Value = "AddX function reference is invalid."
$ws.Range("D37").Value = "UnscentedKalmanFilter_New_FuncGroup.vi"

# Row 38 (keeps the wrap style already on C38)
$ws.Range("C38").Value = "Trouble creating new Diff Drive Pose Est"
$ws.Range("D38").Value = "DiffDrivePoseEst_New.vi"

# Row 39 (distinct font/style -> set WrapText explicitly which creates the new style)
$ws.Range("C39").WrapText = $true
$ws.Range("C39").Value = "Trouble creating new Swerve Drive Pose Est"
$ws.Range("D39").Value = "SwerveDrivePoseEst_New.vi"

# Row 40
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C40").PasteSpecial(-4122) | Out-Null
$ws.Range("C40").Value = "H function reference is invalid."
$ws.Range("D40").Value = "UnscentedKalmanFIlter_Correct_FuncGroup.vi"

# Row 41
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C41").PasteSpecial(-4122) | Out-Null
$ws.Range("C41").Value = "MeanY function reference is invalid."
$ws.Range("D41").Value = "UnscentedKalmanFIlter_Correct_FuncGroup.vi"

# Row 42
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C42").PasteSpecial(-4122) | Out-Null
$ws.Range("C42").Value = "ResidY function reference is invalid."
$ws.Range("D42").Value = "UnscentedKalmanFIlter_Correct_FuncGroup.vi"

# Row 43
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C43").PasteSpecial(-4122) | Out-Null
$ws.Range("C43").Value = "ResidX function reference is invalid."
$ws.Range("D43").Value = "UnscentedKalmanFIlter_Correct_FuncGroup.vi"

# Row 44: also needs the B formula since previously it had no B value at all
$ws.Range("B44").Formula = "=B43+1"
$ws.Range("D7").Copy() | Out-Null
$ws.Range("C44").PasteSpecial(-4122) | Out-Null
$ws.Range("C44").Value = "AddX function reference is invalid."
$ws.Range("D44").Value = "UnscentedKalmanFIlter_Correct_FuncGroup.vi"

# Row 45
$ws.Range("B45").Formula = "=B44+1"
$ws.Range("C45").Value = "F function reference is invalid"
$ws.Range("D45").Value = "KalmanFilterLatencyComp_ApplyPastGlobalMeas_FuncGroup.vi"

# Row 46
$ws.Range("B46").Formula = "=B45+1"
$ws.Range("C46").Value = "H function reference is invalid."
$ws.Range("D46").Value = "KalmanFilterLatencyComp_ApplyPastGlobalMeas_FuncGroup.vi"

# Row 47
$ws.Range("B47").Formula = "=B46+1"
$ws.Range("C47").Value = "VisionCorrect function reference is invalid."
$ws.Range("D47").Value = "KalmanFilterLatencyComp_ApplyPastGlobalMeas_FuncGroup.vi"

# ---------------------------------------------------------------------------
# Rows 48-65: just extend the running Code (B) formula/number series; the
# Description/Originates-from/Note columns remain blank.
# ---------------------------------------------------------------------------
for ($r = 48; $r -le 65; $r++) {
    $prev = $r - 1
    $ws.Range("B$r").Formula = "=B$prev+1"
}

# ---------------------------------------------------------------------------
# Move the active selection to match where editing ended up (C47).
# ---------------------------------------------------------------------------
$ws.Range("C47").Select() | Out-Null
